# Order of writes matters for shared-string index allocation; match
# the order the strings first appear in xl/sharedStrings.xml:
# Caps(A3), TL431(A5), UC3842(A6), Resistor(A4), 30V Zener(A7), 15V Zener(A8), 2N2222 x2(A9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set A3: "Caps" underlined + rest normal (rich text), wrap text
$ws.Range("A3").Value = "Caps`n- 1u x3 (25V for input x1)`n-1n x3`n-2n x1`n-330n(50V) x1`n-330p(150V) x1`n-47u or 68u (63V for output) x1"
$ws.Range("A3").Characters(1, 4).Font.Underline = $true
$ws.Range("A3").Characters(1, 4).Font.Name = "Calibri"
$ws.Range("A3").Characters(1, 4).Font.Size = 11
$ws.Range("A3").Characters(5, 102).Font.Name = "Calibri"
$ws.Range("A3").Characters(5, 102).Font.Size = 11
$ws.Range("A3").WrapText = $true
$ws.Range("A3").RowHeight = 100.8

# Set simple text values
$ws.Range("A5").Value = "TL431"
$ws.Range("A6").Value = "UC3842 "

# Set A4: "Resistor" underlined + rest normal (rich text), wrap text
$ws.Range("A4").Value = "Resistor`n-480R (min 2W) x1`n-0.05R (5W sense)x1`n-10k x5 (x1 %1 tolerance for controller)`n-1k x1`n-600k x1`n-182k x1 (%1 tolerance)`n-14.3k x1`n-22.2k x2`n-10R x1`n-100R x1`n-4.7k x1`n-5k x1"
$ws.Range("A4").Characters(1, 8).Font.Underline = $true
$ws.Range("A4").Characters(1, 8).Font.Name = "Calibri"
$ws.Range("A4").Characters(1, 8).Font.Size = 11
$ws.Range("A4").Characters(9, 172).Font.Name = "Calibri"
$ws.Range("A4").Characters(9, 172).Font.Size = 11
$ws.Range("A4").WrapText = $true
$ws.Range("A4").RowHeight = 187.2

$ws.Range("A7").Value = "30V Zener"
$ws.Range("A8").Value = "15V Zener (Depending on max switch voltage)"
$ws.Range("A9").Value = "2N2222 x2"

# Column C width change (~70.8 chars, matches "best fit" width for the
# long multi-line description text now stored in column A's cells)
$ws.Columns("C").ColumnWidth = 70

# sheet view changes
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A12").Select()
